$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 572.45654
$ws.Range("J17").Value = 377.1842
$ws.Range("L17").Value = 1131.5526
$ws.Range("N17").Value = -1467.5526

$ws.Range("H100").Value = 50001370
$ws.Range("I100").Value = 50001370
$ws.Range("K100").Value = 50001370
$ws.Range("M100").Value = -50000829

$ws.Range("H112").Value = 1255.7966
$ws.Range("J112").Value = 1301.6428
$ws.Range("L112").Value = 3904.9284
$ws.Range("N112").Value = -6120.928400000001

$ws.Range("H115").Value = 1178.75
$ws.Range("I115").Value = 1178.75
$ws.Range("K115").Value = 3536.25
$ws.Range("M115").Value = -1969.25

$ws.Range("H132").Value = 28576270
$ws.Range("I132").Value = 32262564
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 96787692
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -96785162
$ws.Range("N132").Value = -27560

$ws.Range("H138").Value = 2544.86
$ws.Range("I138").Value = 729.9286
$ws.Range("J138").Value = 2840.314
$ws.Range("K138").Value = 2189.7858
$ws.Range("L138").Value = 8520.941999999999
$ws.Range("M138").Value = 2950.2142
$ws.Range("N138").Value = -18800.942

$ws.Range("H141").Value = 155652.47
$ws.Range("I141").Value = 168223.5
$ws.Range("K141").Value = 504670.5
$ws.Range("M141").Value = -499490.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3176.6
$ws.Range("I45").Value = 3625.75
$ws.Range("J45").Value = 1380
$ws.Range("K45").Value = 3625.75
$ws.Range("L45").Value = 1380
$ws.Range("M45").Value = -3248.75
$ws.Range("N45").Value = -2134

$ws.Range("H74").Value = 4908.52
$ws.Range("I74").Value = 6227.933
$ws.Range("J74").Value = 2929.4
$ws.Range("K74").Value = 6227.933
$ws.Range("L74").Value = 2929.4
$ws.Range("M74").Value = -5353.933
$ws.Range("N74").Value = -4677.4

$ws.Range("H77").Value = 4908.52
$ws.Range("I77").Value = 6227.933
$ws.Range("J77").Value = 2929.4
$ws.Range("K77").Value = 31139.665
$ws.Range("L77").Value = 14647
$ws.Range("M77").Value = -26771.665
$ws.Range("N77").Value = -23383

$ws.Range("H132").Value = 4204.25
$ws.Range("I132").Value = 1970.3334
$ws.Range("K132").Value = 5911.0002
$ws.Range("M132").Value = -3381.0002

$ws.Range("H139").Value = 43851.152
$ws.Range("J139").Value = 43851.152
$ws.Range("L139").Value = 43851.152
$ws.Range("N139").Value = -54131.152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22222926
$ws.Range("I16").Value = 27778534
$ws.Range("K16").Value = 27778534
$ws.Range("M16").Value = -27778247

$ws.Range("H31").Value = 2663.5518
$ws.Range("I31").Value = 978.0952
$ws.Range("J31").Value = 7087.875
$ws.Range("K31").Value = 978.0952
$ws.Range("L31").Value = 7087.875
$ws.Range("M31").Value = -683.0952
$ws.Range("N31").Value = -7677.875

$ws.Range("H34").Value = 2663.5518
$ws.Range("I34").Value = 978.0952
$ws.Range("J34").Value = 7087.875
$ws.Range("K34").Value = 978.0952
$ws.Range("L34").Value = 7087.875
$ws.Range("M34").Value = -776.0952
$ws.Range("N34").Value = -7491.875

$ws.Range("H58").Value = 2752.318
$ws.Range("I58").Value = 1633.4364
$ws.Range("J58").Value = 8346.727999999999
$ws.Range("K58").Value = 1633.4364
$ws.Range("L58").Value = 8346.727999999999
$ws.Range("M58").Value = -1430.4364
$ws.Range("N58").Value = -8752.727999999999

$ws.Range("H105").Value = 1489.6666
$ws.Range("I105").Value = 1217.5416
$ws.Range("K105").Value = 1217.5416
$ws.Range("M105").Value = 529.4584

$ws.Range("H113").Value = 22222926
$ws.Range("I113").Value = 27778534
$ws.Range("K113").Value = 27778534
$ws.Range("M113").Value = -27776364

$ws.Range("H132").Value = 3018.7058
$ws.Range("I132").Value = 1755.3846
$ws.Range("K132").Value = 5266.1538
$ws.Range("M132").Value = -2736.1538

$ws.Range("H134").Value = 5863.72
$ws.Range("I134").Value = 6893.8237
$ws.Range("K134").Value = 20681.4711
$ws.Range("M134").Value = -18146.4711

$ws.Range("H136").Value = 2752.318
$ws.Range("I136").Value = 1633.4364
$ws.Range("J136").Value = 8346.727999999999
$ws.Range("K136").Value = 4900.3092
$ws.Range("L136").Value = 25040.184
$ws.Range("M136").Value = -2350.3092
$ws.Range("N136").Value = -30140.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23723.076
$ws.Range("J4").Value = 900
$ws.Range("L4").Value = 2700
$ws.Range("N4").Value = -2924

$ws.Range("H34").Value = 26775.334
$ws.Range("I34").Value = 27604
$ws.Range("J34").Value = 26361
$ws.Range("K34").Value = 82812
$ws.Range("L34").Value = 79083
$ws.Range("M34").Value = -82728
$ws.Range("N34").Value = -79251

$ws.Range("H39").Value = 12242.167
$ws.Range("J39").Value = 12491.706
$ws.Range("L39").Value = 37475.118
$ws.Range("N39").Value = -38063.118

$ws.Range("H55").Value = 4386.5386
$ws.Range("I55").Value = 816.6667
$ws.Range("J55").Value = 5457.5
$ws.Range("K55").Value = 2450.0001
$ws.Range("L55").Value = 16372.5
$ws.Range("M55").Value = -2273.0001
$ws.Range("N55").Value = -16726.5

$ws.Range("H106").Value = 3621.5386
$ws.Range("J106").Value = 3621.5386
$ws.Range("L106").Value = 10864.6158
$ws.Range("N106").Value = -12756.6158

$ws.Range("H129").Value = 2663.6667
$ws.Range("J129").Value = 2057.3845
$ws.Range("L129").Value = 6172.1535
$ws.Range("N129").Value = -16172.1535

$ws.Range("H140").Value = 3492.7856
$ws.Range("I140").Value = 3492.7856
$ws.Range("K140").Value = 10478.3568
$ws.Range("M140").Value = -5298.356800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3278.21
$ws.Range("I126").Value = 2952.8647
$ws.Range("J126").Value = 4204.1924
$ws.Range("K126").Value = 8858.5941
$ws.Range("L126").Value = 12612.5772
$ws.Range("M126").Value = -6388.5941
$ws.Range("N126").Value = -17552.5772

$ws.Range("H135").Value = 38027.5
$ws.Range("J135").Value = 38027.5
$ws.Range("L135").Value = 38027.5
$ws.Range("N135").Value = -48167.5

$ws.Range("H140").Value = 38991.875
$ws.Range("J140").Value = 38991.875
$ws.Range("L140").Value = 38991.875
$ws.Range("N140").Value = -49351.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 27499.5
$ws.Range("J26").Value = 27499.5
$ws.Range("L26").Value = 27499.5
$ws.Range("N26").Value = -28089.5

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

$ws.Range("H138").Value = 51083.8
$ws.Range("J138").Value = 51083.8
$ws.Range("L138").Value = 51083.8
$ws.Range("N138").Value = -61363.8

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

$ws.Range("H141").Value = 36071.43
$ws.Range("J141").Value = 36071.43
$ws.Range("L141").Value = 36071.43
$ws.Range("N141").Value = -46431.43

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3179.081
$ws.Range("I136").Value = 1557.3667
$ws.Range("J136").Value = 10129.286
$ws.Range("K136").Value = 4672.1001
$ws.Range("L136").Value = 30387.858
$ws.Range("M136").Value = -2122.1001
$ws.Range("N136").Value = -35487.858

$ws.Range("H138").Value = 38549.668
$ws.Range("J138").Value = 38549.668
$ws.Range("L138").Value = 38549.668
$ws.Range("N138").Value = -48829.668
